# Generate Report for Handback
#
# - Update the "Ready for handoff" status everywhere to
#   "Handed back: in sync with en-US".
# - Fill in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns on the zh-cn and de-de sheets,
#   now that the handback round-trip produced real xlf files + timestamps,
#   and turn "Latest Target File" into a hyperlink to the source doc
#   (mirrors column A).
# - Widen the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

# 1) Status text changed workbook-wide: "Ready for handoff" -> "Handed back: in sync with en-US"
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# 2) zh-cn sheet: fill Latest Target File / Latest Handback File / Latest Handback DateTime
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("J2").Value = "3b080f16-260d-40c7-bb15-5300fe05dacc.b96a7722e32f43d36ae1594a91fb3d3bf74bb8c4.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-19 16:34:39"

$zh.Range("J3").Value = "e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.8af04b4c67b9327e047b15a9b570d66be003b48d.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-19 16:34:39"

# Rebuild the hyperlinks collection in the desired order: A2, I2, A3, I3
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/3b080f16-260d-40c7-bb15-5300fe05dacc.md", "", "", "3b080f16-260d-40c7-bb15-5300fe05dacc.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/3b080f16-260d-40c7-bb15-5300fe05dacc.md", "", "", "3b080f16-260d-40c7-bb15-5300fe05dacc.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md", "", "", "e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md", "", "", "e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md")

# Widen columns C (Status), I (Latest Target File), J (Latest Handback File)
$zh.Range("C1").ColumnWidth = 29.16
$zh.Range("I1").ColumnWidth = 39.16
$zh.Range("J1").ColumnWidth = 39.16

# 3) de-de sheet: fill Latest Target File / Latest Handback File / Latest Handback DateTime
$de = $wb.Worksheets.Item("de-de")

$de.Range("J2").Value = "3b080f16-260d-40c7-bb15-5300fe05dacc.b96a7722e32f43d36ae1594a91fb3d3bf74bb8c4.de-de.xlf"
$de.Range("K2").Value = "2016-08-19 16:34:45"

$de.Range("J3").Value = "e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.8af04b4c67b9327e047b15a9b570d66be003b48d.de-de.xlf"
$de.Range("K3").Value = "2016-08-19 16:34:45"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/3b080f16-260d-40c7-bb15-5300fe05dacc.md", "", "", "3b080f16-260d-40c7-bb15-5300fe05dacc.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/3b080f16-260d-40c7-bb15-5300fe05dacc.md", "", "", "3b080f16-260d-40c7-bb15-5300fe05dacc.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md", "", "", "e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/efa98428442028819f4de428d12d72932a522559/e2e/e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md", "", "", "e34b2465-6c77-4eaa-b6a4-6aa9d0fb5a7b.md")

$de.Range("C1").ColumnWidth = 29.16
$de.Range("I1").ColumnWidth = 39.16
$de.Range("J1").ColumnWidth = 39.16

# 4) Overview sheet: widen the zh-cn / de-de status columns (E, F) to match the new, longer status text
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E1").ColumnWidth = 29.16
$ov.Range("F1").ColumnWidth = 29.16
